$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nam")

# Fix the mobile number / name typo for "Nam Tich"
$ws.Range("C463").Value = "Nam Tịch"
$ws.Range("F463").Value = "'0818592855"

# Remove the "Ty Long An" row entirely (row 600)
$ws.Rows.Item(600).Delete()
